# "replaced direct usage of inverse coefficients"
#
# Row 83 used to hold a single work interval 12:30-13:15. That interval is
# split into two separate intervals (a short break is carved out of the
# middle): 12:30-13:00 (stays on row 83) and 13:15-14:30 (new row 84).
# Everything below (the blank separator row and the three summary rows)
# shifts down by one row, and the running SUM() is extended to pick up the
# newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 84, pushing the old rows 84-87 down to 85-88.
# Dependent formulas (e.g. the SUM in what is now F86) are adjusted
# automatically by the engine to keep referencing the "last data row".
$ws.Rows(84).Insert()

# Row 83: shorten the first interval's end time 13:15 -> 13:00.
$ws.Range("E83").Value = 0.54166666666666663

# Row 84: brand-new data row for the carved-out second interval
# (13:15 -> 14:30), same date as row 83.
$ws.Range("A84").Value = 2014
$ws.Range("B84").Value = 3
$ws.Range("C84").Value = 21
$ws.Range("D84").Value = 0.55208333333333337
$ws.Range("E84").Value = 0.60416666666666663

# Duration formulas for the new row, matching the shared pattern used by
# the rest of the F/G columns.
$ws.Range("F84").Formula = "=(E84-D84)*24*60"
$ws.Range("G84").Formula = "=F84/60"

# Keep the persisted selection in sync with the diff (it tracked row 84,
# now row 85 after the insert).
$ws.Range("A85").Select()
